$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("target number of stations" and
# everything to its right shifts one column to the right).
$ws.Columns("E:E").Insert()

# New column E header + values: "station annual capacity factor"
$ws.Cells.Item(1, 5).Value = "station annual capacity factor"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(4, 5).Value = 1

# New scenario row 5: "closed loop 80 cap factor" - same inputs as row 4
# ("closed loop") but with the station annual capacity factor set to 0.8.
$ws.Cells.Item(5, 1).Value = "closed loop 80 cap factor"
$ws.Cells.Item(5, 2).Formula = "=B4+1"
$ws.Cells.Item(5, 3).Value = 2022
$ws.Cells.Item(5, 4).Value = 1000
$ws.Cells.Item(5, 5).Value = 0.8
$ws.Cells.Item(5, 6).Value = 10
$ws.Cells.Item(5, 7).Value = 100
$ws.Cells.Item(5, 8).Formula = "=17.09/100"
$ws.Cells.Item(5, 9).Value = 6.0279999999999996
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 10.18
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0.31
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = "baseline"
$ws.Cells.Item(5, 17).Value = 366.15
$ws.Cells.Item(5, 18).Value = 105
$ws.Cells.Item(5, 19).Value = 1
$ws.Cells.Item(5, 20).Value = 1
$ws.Cells.Item(5, 21).Value = 1
$ws.Cells.Item(5, 22).Value = 53
$ws.Cells.Item(5, 23).Value = 5450
$ws.Cells.Item(5, 24).Value = 1
$ws.Cells.Item(5, 25).Value = 0
$ws.Cells.Item(5, 26).Value = 0
$ws.Cells.Item(5, 27).Value = 5250
$ws.Cells.Item(5, 28).Value = 0.25
$ws.Cells.Item(5, 29).Value = 300
$ws.Cells.Item(5, 30).Value = 1
$ws.Cells.Item(5, 31).Value = 0.99990000000000001
$ws.Cells.Item(5, 32).Value = 0.0731028611028611
$ws.Cells.Item(5, 33).Value = 1
$ws.Cells.Item(5, 34).Value = 9.6467120334224301
$ws.Cells.Item(5, 35).Value = 3500
$ws.Cells.Item(5, 36).Value = 1
$ws.Cells.Item(5, 37).Value = 0
$ws.Cells.Item(5, 38).Value = 0
$ws.Cells.Item(5, 39).Value = 1

$ws.Range("A6").Select()
